# Updated cryptos list on Sun May 21 13:49:58 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# table with a newer snapshot, including a couple of rows (40/41) whose
# rank order flipped (VeChain <-> FraxShare).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # All "Price" cells are stored as literal text (e.g. "0.3650",
    # "27.076.62", "0.000008718"). Excel's normal Range.Value setter
    # auto-detects numeric-looking strings and coerces them into real
    # numbers, which silently drops meaningful trailing zeros and
    # re-renders multi-dot "thousands" separated strings. Force the
    # cell to Text, write the literal string, then restore the
    # original style so no stray number formatting is left behind.
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

# row -> (Price, Volume(1h)) updates
$priceVolume = @{
    2  = @("27.076.62",     "  +0.35%  ")
    3  = @("1.824.50",      "  +0.17%  ")
    4  = @("1.007",         "  +0.37%  ")
    5  = @("312.58",        "  +0.47%  ")
    7  = @("0.4686",        "  +0.15%  ")
    8  = @("0.3650",        "  -0.43%  ")
    9  = @("0.07381",       "  +0.35%  ")
    10 = @("0.8780",        "  +0.50%  ")
    11 = @("20.24",         "  -0.30%  ")
    12 = @("1.895.98",      "  +4.21%  ")
    13 = @("0.07523",       "  +5.24%  ")
    14 = @("5.368",         "  -1.02%  ")
    15 = @("92.76",         "  +1.24%  ")
    16 = @("6.522",         "  +0.08%  ")
    20 = @("27.450.49",     "  +1.65%  ")
    22 = @("5.231",         "  -1.19%  ")
    24 = @("2.076.10",      "  +1.60%  ")
    25 = @("1.881",         "  -0.59%  ")
    26 = @("151.42",        "  +0.29%  ")
    27 = @("18.49",         "  +0.40%  ")
    28 = @("2.131",         "  -0.50%  ")
    30 = @("116.37",        "  -0.27%  ")
    34 = @("4.505",         "  +0.05%  ")
    37 = @("2.517",         "  +5.82%  ")
    38 = @("1.089",         "  -0.49%  ")
    39 = @("0.05293",       "  -0.34%  ")
    42 = @("2.933",         "  -1.51%  ")
    43 = @("0.5242",        "  -1.12%  ")
    45 = @("8.361",         "  -1.27%  ")
    46 = @("0.4889",        "  -0.17%  ")
    47 = @("10.41",         "  -0.79%  ")
    49 = @("104.17",        "  +1.05%  ")
    51 = @("0.06264",       "  -0.51%  ")
}

# rows where only Price changed (Volume(1h) untouched)
$priceOnly = @{
    18 = "0.000008718"
    31 = "0.08918"
}

# rows where only Volume(1h) changed (Price untouched)
$volumeOnly = @{
    6  = "  +0.37%  "
    17 = "  +0.15%  "
    21 = "  -0.52%  "
    23 = "  +0.01%  "
    29 = "  -1.53%  "
    32 = "  -1.42%  "
    33 = "  +0.01%  "
    35 = "  -0.07%  "
    36 = "  +0.42%  "
    44 = "  -0.80%  "
    48 = "  +0.40%  "
}

foreach ($row in $priceVolume.Keys) {
    $vals = $priceVolume[$row]
    Set-TextValue $ws.Range("D$row") $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

foreach ($row in $priceOnly.Keys) {
    Set-TextValue $ws.Range("D$row") $priceOnly[$row]
}

foreach ($row in $volumeOnly.Keys) {
    $ws.Range("E$row").Value = $volumeOnly[$row]
}

# Rows 40/41 swapped places in the source ranking (row 40 is now VeChain,
# row 41 is now FraxShare) in addition to their price/volume refresh.
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.01930"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "7.313"
$ws.Range("E41").Value = "  +1.74%  "

# Row 50: only Price changed.
Set-TextValue $ws.Range("D50") "1.648"
